{"js": "// Replace the 25 \"three-digit \u00f7 one-digit\" practice answers in the table\n// with the regenerated values, keeping everything else untouched.\nconst replacements = [\n  [\"118\u00f75=23, 3\", \"967\u00f78=120, 7\"],\n  [\"859\u00f76=143, 1\", \"914\u00f77=130, 4\"],\n  [\"195\u00f76=32, 3\", \"196\u00f76=32, 4\"],\n  [\"307\u00f74=76, 3\", \"257\u00f73=85, 2\"],\n  [\"942\u00f74=235, 2\", \"853\u00f74=213, 1\"],\n  [\"782\u00f72=391, 0\", \"467\u00f73=155, 2\"],\n  [\"218\u00f74=54, 2\", \"714\u00f78=89, 2\"],\n  [\"541\u00f79=60, 1\", \"999\u00f76=166, 3\"],\n  [\"153\u00f78=19, 1\", \"210\u00f73=70, 0\"],\n  [\"684\u00f77=97, 5\", \"958\u00f77=136, 6\"],\n  [\"215\u00f79=23, 8\", \"752\u00f72=376, 0\"],\n  [\"314\u00f75=62, 4\", \"681\u00f73=227, 0\"],\n  [\"379\u00f79=42, 1\", \"366\u00f73=122, 0\"],\n  [\"834\u00f72=417, 0\", \"869\u00f73=289, 2\"],\n  [\"445\u00f76=74, 1\", \"153\u00f79=17, 0\"],\n  [\"449\u00f75=89, 4\", \"117\u00f72=58, 1\"],\n  [\"251\u00f79=27, 8\", \"624\u00f74=156, 0\"],\n  [\"834\u00f76=139, 0\", \"394\u00f72=197, 0\"],\n  [\"425\u00f73=141, 2\", \"503\u00f79=55, 8\"],\n  [\"102\u00f72=51, 0\", \"669\u00f74=167, 1\"],\n  [\"884\u00f77=126, 2\", \"597\u00f77=85, 2\"],\n  [\"126\u00f77=18, 0\", \"540\u00f72=270, 0\"],\n  [\"689\u00f72=344, 1\", \"233\u00f73=77, 2\"],\n  [\"887\u00f77=126, 5\", \"650\u00f76=108, 2\"],\n  [\"842\u00f76=140, 2\", \"599\u00f74=149, 3\"],\n];\n\n// Use search+Range.insertText(\"Replace\") so the run's existing formatting\n// (font, size, etc.) is preserved instead of being reset by a clear+insert.\nfor (const [before, after] of replacements) {\n  const results = context.document.body.search(before, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(after, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 \"three-digit \u00f7 one-digit\" practice answers in the table\n# with the regenerated values, keeping everything else (fonts, sizes,\n# alignment) untouched by using Find/Replace on each table cell's Range.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"118\u00f75=23, 3\", \"967\u00f78=120, 7\"),\n    @(\"859\u00f76=143, 1\", \"914\u00f77=130, 4\"),\n    @(\"195\u00f76=32, 3\", \"196\u00f76=32, 4\"),\n    @(\"307\u00f74=76, 3\", \"257\u00f73=85, 2\"),\n    @(\"942\u00f74=235, 2\", \"853\u00f74=213, 1\"),\n    @(\"782\u00f72=391, 0\", \"467\u00f73=155, 2\"),\n    @(\"218\u00f74=54, 2\", \"714\u00f78=89, 2\"),\n    @(\"541\u00f79=60, 1\", \"999\u00f76=166, 3\"),\n    @(\"153\u00f78=19, 1\", \"210\u00f73=70, 0\"),\n    @(\"684\u00f77=97, 5\", \"958\u00f77=136, 6\"),\n    @(\"215\u00f79=23, 8\", \"752\u00f72=376, 0\"),\n    @(\"314\u00f75=62, 4\", \"681\u00f73=227, 0\"),\n    @(\"379\u00f79=42, 1\", \"366\u00f73=122, 0\"),\n    @(\"834\u00f72=417, 0\", \"869\u00f73=289, 2\"),\n    @(\"445\u00f76=74, 1\", \"153\u00f79=17, 0\"),\n    @(\"449\u00f75=89, 4\", \"117\u00f72=58, 1\"),\n    @(\"251\u00f79=27, 8\", \"624\u00f74=156, 0\"),\n    @(\"834\u00f76=139, 0\", \"394\u00f72=197, 0\"),\n    @(\"425\u00f73=141, 2\", \"503\u00f79=55, 8\"),\n    @(\"102\u00f72=51, 0\", \"669\u00f74=167, 1\"),\n    @(\"884\u00f77=126, 2\", \"597\u00f77=85, 2\"),\n    @(\"126\u00f77=18, 0\", \"540\u00f72=270, 0\"),\n    @(\"689\u00f72=344, 1\", \"233\u00f73=77, 2\"),\n    @(\"887\u00f77=126, 5\", \"650\u00f76=108, 2\"),\n    @(\"842\u00f76=140, 2\", \"599\u00f74=149, 3\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    # wdFindContinue=1, wdReplaceAll=2 -- scan the whole story, swap every\n    # exact (case-sensitive) occurrence, leave the run's own formatting as-is.\n    $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n"}
